$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 260-261; this pushes the existing rows
# 260-303 down to 262-305 (and widens the used range to A1:R305).
$ws.Rows("260:261").Insert()

# Populate the two newly inserted rows with the new weekly records.

# Row 260 - Región de O'Higgins
$ws.Cells.Item(260, 1).Value = 10
$ws.Cells.Item(260, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value = "La Araucanía"
$ws.Cells.Item(260, 4).Value = 44522
$ws.Cells.Item(260, 5).Value = 9
$ws.Cells.Item(260, 6).Value = 100112023
$ws.Cells.Item(260, 7).Value = "Brócoli"
$ws.Cells.Item(260, 8).Value = "Sin especificar"
$ws.Cells.Item(260, 9).Value = "Primera"
$ws.Cells.Item(260, 10).Value = 1200
$ws.Cells.Item(260, 11).Value = 800
$ws.Cells.Item(260, 12).Value = 800
$ws.Cells.Item(260, 13).Value = 800
$ws.Cells.Item(260, 14).Value = "$/unidad"
$ws.Cells.Item(260, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(260, 16).Value = 800
$ws.Cells.Item(260, 17).Value = 1
$ws.Cells.Item(260, 18).Value = "Hortaliza"

# Row 261 - Región del Maule
$ws.Cells.Item(261, 1).Value = 10
$ws.Cells.Item(261, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(261, 3).Value = "La Araucanía"
$ws.Cells.Item(261, 4).Value = 44522
$ws.Cells.Item(261, 5).Value = 9
$ws.Cells.Item(261, 6).Value = 100112023
$ws.Cells.Item(261, 7).Value = "Brócoli"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 2000
$ws.Cells.Item(261, 11).Value = 800
$ws.Cells.Item(261, 12).Value = 800
$ws.Cells.Item(261, 13).Value = 800
$ws.Cells.Item(261, 14).Value = "$/unidad"
$ws.Cells.Item(261, 15).Value = "Región del Maule"
$ws.Cells.Item(261, 16).Value = 800
$ws.Cells.Item(261, 17).Value = 1
$ws.Cells.Item(261, 18).Value = "Hortaliza"
